$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 'AMD Radeon Rx 6800xt 16G Pcie 4.0'
$ws.Range("C2").Value = 'https://www.ebay.com/itm/175531714780?hash=item28de8168dc:g:t0YAAOSwBbpjmF4E&amdata=enc%3AAQAHAAAAoP4qtSmdnF7K6kKgCTQO2HJR52%2BBOgy%2BA%2FQHGLUkOvnoonZ9kbd%2BT08dcMptHax7gC%2BprXsM6AAI8ii5YMD7tw6n%2FtO%2FcMfxF%2Blcba81qJg4ifs1Eek5fSoq%2F3f6XNj08jtuJbQipNh1AHfxu5BFoKqxgqE2pyZufALb2N86QBX7ndCGJlHfKCs837lFkH1Dzh2RosD4Z%2FZvq8UayNlgw4U%3D%7Ctkp%3ABk9SR6iqv8WjYQ'
$ws.Range("E2").Value = 565
$ws.Range("H2").Value = 590

# Row 3 changes
$ws.Range("C3").Value = 'https://www.ebay.com/itm/385163444573?epid=23050214728&hash=item59ad87c15d:g:oFAAAOSw5yxjRIsG&amdata=enc%3AAQAHAAAAoHHTmOg%2BoR3LtRCMXYQ8lAEEmjWes9DhmmiUWq2HMpjVj7cHG6OArdLFELfIpHbYXf%2FxvqhhQHI16MdjiGX6kObyLNflfSBpAixGObFlZPIjcOfiKrnYEzcbQ3NhsODrompGmh%2F9TnPL6rVWiH%2F4z%2BGgeGZn23lV0AsL3%2Bn%2BHb5b5%2BmGMfvoAiHHIZ3wCCWloN1MM%2B9v071%2BjQBFDA10B0s%3D%7Ctkp%3ABk9SR6iqv8WjYQ'

# Row 4 changes
$ws.Range("C4").Value = 'https://www.ebay.com/itm/364076525731?hash=item54c4a6e8a3:g:JX4AAOSw~CJjlfzi&amdata=enc%3AAQAHAAAA4NEFSvBjjHf1HqZEVcXoZ9aEHAWu9bZtDogcLQ1aECJrJ4bqtAM2GcxfCgyJWBDbzJoIrzPMG0Cg7y%2BGOVailhMwXU4A7A09ESoKAM2JBOuXId253o%2F8uJGTEwKxMubxgRRQ33B60oIIfIMx3Ec2hn5I4WNrnkMOXgGZwWTZiRioF2abPA%2F%2BoOko23se95%2F6Oh8xhM6O4mxI7a3Tsxq9hVnxbxT%2FSHjOsudJm8D6BwbvIrMmxij2Ty0e4UznU2XJQydJRsov8IRvPYjb1HcWArAdo3YfvIE%2Fl4dg4iNfz88k%7Ctkp%3ABk9SR6iqv8WjYQ'
$ws.Range("E4").Value = 454
$ws.Range("H4").Value = 469
$ws.Range("K4").Value = '1d 20h restantes'

# Row 5 changes
$ws.Range("A5").Value = 1
$ws.Range("C5").Value = 'https://www.ebay.com/itm/155308020238?hash=item242914520e:g:HecAAOSw6CBjmLE8&amdata=enc%3AAQAHAAAAkINXBQ6TDT7aR4w2cDXsHla8%2F9Ul%2FpYukvcoYW7Eov4Cb4GvJmWTkx5qmIHfV0oZdazVqMWw23GMNqSOI35n6NV1iZnDB6D4YsX3gaVZKEuilkg303CshW7W6bHiny8HRZQa1EsmTT3oBK3AlIB%2BW%2Bh2KYyVY4mCHZO7I4MWGdZUBzVdje8cIy%2BwlEea%2BHuJTg%3D%3D%7Ctkp%3ABk9SR6iqv8WjYQ'
$ws.Range("K5").Value = '1d 22h restantes'

# Row 6 changes
$ws.Range("C6").Value = 'https://www.ebay.com/itm/175531638048?hash=item28de803d20:g:1rwAAOSwGZ5jjxIj&amdata=enc%3AAQAHAAAAkDUDHPwW2TkeJpPu0PLgguIKqM6OK8we5QyqR8AY203rotyF%2Bfy4InheLSYDL5LLc9d3DYd3LHUgqYgGK8IGy5vKJ4mIijbeK0YIA426rAE2%2B626WktOUh3QdtVk5rLxiGE4lxeCUYdWxX%2F6ttQ3jyhHQWpmafQEfEuGL19D%2BfeY3etONS5bI3Epg5BMONOPWA%3D%3D%7Ctkp%3ABk9SR6iqv8WjYQ'
$ws.Range("K6").Value = '3d 15h restantes'

